$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CSC-CSCSoCECBiaSY")
$ws.Activate()

# Update row 8 values from 0.33 to 0.25 for columns B through AE (solar PV es share of existing capacity)
$ws.Range("B8:AE8").Value = 0.25

# Update the view state: scroll position and selection
$ws.Application.ActiveWindow.ScrollColumn = 6
$ws.Range("AF8:AG8").Select()
